$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted above the existing row 104,
# pushing every subsequent record (old rows 104-229) down by one row
# (new rows 105-230). Insert a fresh row at position 104 to reproduce
# that shift, then populate it with the new record's data.
$ws.Rows(104).Insert()

$ws.Range('A104').Value = 11
$ws.Range('B104').Value = 'Vega Monumental Concepción'
$ws.Range('C104').Value = 'Bíobío'
$ws.Range('D104').Value = 44895
$ws.Range('E104').Value = 8
$ws.Range('F104').Value = 100112003
$ws.Range('G104').Value = 'Ajo'
$ws.Range('H104').Value = 'Chino'
$ws.Range('I104').Value = 'Primera'
$ws.Range('J104').Value = 500
$ws.Range('K104').Value = 14000
$ws.Range('L104').Value = 15000
$ws.Range('M104').Value = 14600
$ws.Range('N104').Value = '$/caja 10 kilos'
$ws.Range('O104').Value = 'China'
$ws.Range('P104').Value = 1460
$ws.Range('Q104').Value = 10
$ws.Range('R104').Value = 'Hortaliza'
